$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CLIENTES table (rows 2-7) ---
# "cedula" type: INTEGER -> BIGINT
$ws.Range("B3").Value = "BIGINT"
# "foto" type: VARCHAR(50) -> INTEGER
$ws.Range("F3").Value = "INTEGER"
# "foto" dominio (domain): (blank) -> "1, 0"
$ws.Range("F6").Value = "1, 0"
# "foto" observaciones: (blank) -> "Booleano"
$ws.Range("F7").Value = "Booleano"

# --- FACTURAS table (rows 10-15) ---
# "codigo" type: INTEGER -> BIGINT
$ws.Range("B11").Value = "BIGINT"
# "cedula_cliente" type: INTEGER -> BIGINT
$ws.Range("C11").Value = "BIGINT"

# --- PRODUCTOS table (rows 18-23) ---
# "codigo" type: INTEGER -> BIGINT
$ws.Range("B19").Value = "BIGINT"

# --- DETALLES table (rows 26-31) ---
# "id" type: INTEGER -> BIGINT
$ws.Range("B27").Value = "BIGINT"
# "codigo_factura" type: INTEGER -> BIGINT
$ws.Range("C27").Value = "BIGINT"
# "codigo_producto" type: INTEGER -> BIGINT
$ws.Range("D27").Value = "BIGINT"

# Update the active cell selection to reflect where the author was last working
$ws.Range("G6").Select()
